# Auto-generated edit script applying the committed diff to Belias_Profits sheets.
# Updates price/profit columns (H-N) for specific Leve rows across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 15000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 15000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 45000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -45320
$ws.Range("H88").Value = 27992.125
$ws.Range("I88").Value = 72032.125
$ws.Range("J88").Value = 5972.125
$ws.Range("K88").Value = 72032.125
$ws.Range("L88").Value = 5972.125
$ws.Range("M88").Value = -71626.125
$ws.Range("N88").Value = -6784.125
$ws.Range("H91").Value = 27992.125
$ws.Range("I91").Value = 72032.125
$ws.Range("J91").Value = 5972.125
$ws.Range("K91").Value = 72032.125
$ws.Range("L91").Value = 5972.125
$ws.Range("M91").Value = -70628.125
$ws.Range("N91").Value = -8780.125
$ws.Range("H108").Value = 31000
$ws.Range("J108").Value = 31000
$ws.Range("L108").Value = 31000
$ws.Range("N108").Value = -38680
$ws.Range("H125").Value = 2276
$ws.Range("H137").Value = 1918173.2
$ws.Range("I137").Value = 2267.923
$ws.Range("J137").Value = 3474846.5
$ws.Range("K137").Value = 6803.768999999999
$ws.Range("L137").Value = 10424539.5
$ws.Range("M137").Value = -4253.768999999999
$ws.Range("N137").Value = -10429639.5
$ws.Range("H138").Value = 2614.1943
$ws.Range("I138").Value = 1666.186
$ws.Range("J138").Value = 4019.862
$ws.Range("K138").Value = 4998.558
$ws.Range("L138").Value = 12059.586
$ws.Range("M138").Value = 141.442
$ws.Range("N138").Value = -22339.586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14509.3125
$ws.Range("I32").Value = 8087.973
$ws.Range("K32").Value = 8087.973
$ws.Range("M32").Value = -7800.973
$ws.Range("H45").Value = 12988043
$ws.Range("I45").Value = 15152550
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 15152550
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -15152173
$ws.Range("N45").Value = -1754
$ws.Range("H97").Value = 907.9583
$ws.Range("I97").Value = 950.5294
$ws.Range("J97").Value = 804.5714
$ws.Range("K97").Value = 950.5294
$ws.Range("L97").Value = 804.5714
$ws.Range("M97").Value = -454.5294
$ws.Range("N97").Value = -1796.5714
$ws.Range("H122").Value = 1576.0667
$ws.Range("I122").Value = 1289.6111
$ws.Range("J122").Value = 2005.75
$ws.Range("K122").Value = 3868.8333
$ws.Range("L122").Value = 6017.25
$ws.Range("M122").Value = -1418.8333
$ws.Range("N122").Value = -10917.25
$ws.Range("H132").Value = 1822.6562
$ws.Range("I132").Value = 1141.5
$ws.Range("J132").Value = 3321.2
$ws.Range("K132").Value = 3424.5
$ws.Range("L132").Value = 9963.599999999999
$ws.Range("M132").Value = -894.5
$ws.Range("N132").Value = -15023.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 33335626
$ws.Range("I20").Value = 52633716
$ws.Range("K20").Value = 52633716
$ws.Range("M20").Value = -52633469
$ws.Range("H86").Value = 1786.875
$ws.Range("I86").Value = 1876.8235
$ws.Range("J86").Value = 1568.4286
$ws.Range("K86").Value = 1876.8235
$ws.Range("L86").Value = 1568.4286
$ws.Range("M86").Value = -753.8235
$ws.Range("N86").Value = -3814.4286
$ws.Range("H89").Value = 1786.875
$ws.Range("I89").Value = 1876.8235
$ws.Range("J89").Value = 1568.4286
$ws.Range("K89").Value = 9384.1175
$ws.Range("L89").Value = 7842.143
$ws.Range("M89").Value = -3768.1175
$ws.Range("N89").Value = -19074.143
$ws.Range("H105").Value = 2474.9792
$ws.Range("I105").Value = 2473.8914
$ws.Range("K105").Value = 2473.8914
$ws.Range("M105").Value = -726.8914

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 51722.91
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 54895.2
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 54895.2
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -65255.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1791.2727
$ws.Range("I140").Value = 670.4
$ws.Range("J140").Value = 13000
$ws.Range("K140").Value = 2011.2
$ws.Range("L140").Value = 39000
$ws.Range("M140").Value = 3168.8
$ws.Range("N140").Value = -49360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 15000.5
$ws.Range("J5").Value = 15000.5
$ws.Range("L5").Value = 15000.5
$ws.Range("N5").Value = -15224.5
$ws.Range("H70").Value = 5506.6875
$ws.Range("I70").Value = 5200.636
$ws.Range("K70").Value = 5200.636
$ws.Range("M70").Value = -4930.636
$ws.Range("H73").Value = 5506.6875
$ws.Range("I73").Value = 5200.636
$ws.Range("K73").Value = 5200.636
$ws.Range("M73").Value = -4264.636
$ws.Range("H80").Value = 2376.1904
$ws.Range("I80").Value = 2412.5
$ws.Range("J80").Value = 2353.8462
$ws.Range("K80").Value = 2412.5
$ws.Range("L80").Value = 2353.8462
$ws.Range("M80").Value = -1414.5
$ws.Range("N80").Value = -4349.8462
$ws.Range("H83").Value = 2376.1904
$ws.Range("I83").Value = 2412.5
$ws.Range("J83").Value = 2353.8462
$ws.Range("K83").Value = 12062.5
$ws.Range("L83").Value = 11769.231
$ws.Range("M83").Value = -7070.5
$ws.Range("N83").Value = -21753.231
$ws.Range("H97").Value = 1348.6072
$ws.Range("I97").Value = 1315.0385
$ws.Range("J97").Value = 1785
$ws.Range("K97").Value = 1315.0385
$ws.Range("L97").Value = 1785
$ws.Range("M97").Value = -819.0385000000001
$ws.Range("N97").Value = -2777
$ws.Range("H126").Value = 37038884
$ws.Range("I126").Value = 66667930
$ws.Range("J126").Value = 2575
$ws.Range("K126").Value = 200003790
$ws.Range("L126").Value = 7725
$ws.Range("M126").Value = -200001320
$ws.Range("N126").Value = -12665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 23421.777
$ws.Range("I68").Value = 26224.5
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 26224.5
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -25475.5
$ws.Range("N68").Value = -2498
$ws.Range("H71").Value = 23421.777
$ws.Range("I71").Value = 26224.5
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 131122.5
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -127378.5
$ws.Range("N71").Value = -12488
